$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1897106109324759
$ws.Range("C2").Value = 0.5466237942122186
$ws.Range("J2").Value = 0.05627009646302251
$ws.Range("P2").Value = 0.135048231511254
$ws.Range("S2").Value = 0.07234726688102894
$ws.Range("B3").Value = 0.008333333333333333
$ws.Range("C3").Value = 0.02777777777777778
$ws.Range("J3").Value = 0.08611111111111111
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.1277777777777778
$ws.Range("J4").Value = 0.1234567901234568
$ws.Range("P4").Value = 0.6049382716049383
$ws.Range("S4").Value = 0.2716049382716049
$ws.Range("B6").Value = 0.07306889352818371
$ws.Range("D6").Value = 0.006263048016701462
$ws.Range("E6").Value = 0.00208768267223382
$ws.Range("F6").Value = 0.06889352818371608
$ws.Range("J6").Value = 0.2964509394572025
$ws.Range("O6").Value = 0.01670146137787056
$ws.Range("Q6").Value = 0.1711899791231733
$ws.Range("R6").Value = 0.07933194154488518
$ws.Range("S6").Value = 0.2860125260960334
$ws.Range("B7").Value = 0.1171875
$ws.Range("D7").Value = 0.01041666666666667
$ws.Range("E7").Value = 0.002604166666666667
$ws.Range("F7").Value = 0.046875
$ws.Range("J7").Value = 0.1666666666666667
$ws.Range("O7").Value = 0.0078125
$ws.Range("Q7").Value = 0.1953125
$ws.Range("R7").Value = 0.1015625
$ws.Range("S7").Value = 0.3515625
$ws.Range("B8").Value = 0.09362934362934362
$ws.Range("D8").Value = 0.02123552123552123
$ws.Range("F8").Value = 0.07625482625482626
$ws.Range("J8").Value = 0.1708494208494208
$ws.Range("O8").Value = 0.01061776061776062
$ws.Range("Q8").Value = 0.1988416988416988
$ws.Range("R8").Value = 0.1090733590733591
$ws.Range("S8").Value = 0.3194980694980695
$ws.Range("B9").Value = 0.09254498714652956
$ws.Range("D9").Value = 0.02313624678663239
$ws.Range("F9").Value = 0.08740359897172237
$ws.Range("J9").Value = 0.1670951156812339
$ws.Range("O9").Value = 0.01799485861182519
$ws.Range("Q9").Value = 0.1568123393316195
$ws.Range("R9").Value = 0.1182519280205656
$ws.Range("S9").Value = 0.3367609254498715
$ws.Range("B10").Value = 0.0926812585499316
$ws.Range("D10").Value = 0.01538987688098495
$ws.Range("E10").Value = 0.0006839945280437756
$ws.Range("F10").Value = 0.06121751025991792
$ws.Range("J10").Value = 0.2523939808481532
$ws.Range("O10").Value = 0.02735978112175103
$ws.Range("Q10").Value = 0.2058823529411765
$ws.Range("R10").Value = 0.07900136798905609
$ws.Range("S10").Value = 0.265389876880985
$ws.Range("G11").Value = 0.1341463414634146
$ws.Range("J11").Value = 0.07926829268292683
$ws.Range("K11").Value = 0.1727642276422764
$ws.Range("L11").Value = 0.6036585365853658
$ws.Range("S11").Value = 0.01016260162601626
$ws.Range("G12").Value = 0.7913907284768212
$ws.Range("J12").Value = 0.1556291390728477
$ws.Range("K12").Value = 0.009933774834437087
$ws.Range("L12").Value = 0.02317880794701987
$ws.Range("S12").Value = 0.01986754966887417
$ws.Range("F13").Value = 0.008849557522123894
$ws.Range("G13").Value = 0.7168141592920354
$ws.Range("J13").Value = 0.2566371681415929
$ws.Range("S13").Value = 0.01769911504424779
$ws.Range("F15").Value = 0.02415458937198068
$ws.Range("H15").Value = 0.2028985507246377
$ws.Range("I15").Value = 0.06038647342995169
$ws.Range("J15").Value = 0.3671497584541063
$ws.Range("K15").Value = 0.05797101449275362
$ws.Range("M15").Value = 0.01690821256038647
$ws.Range("N15").Value = 0.002415458937198068
$ws.Range("O15").Value = 0.05797101449275362
$ws.Range("S15").Value = 0.2101449275362319
$ws.Range("F16").Value = 0.01526717557251908
$ws.Range("H16").Value = 0.1908396946564886
$ws.Range("I16").Value = 0.09414758269720101
$ws.Range("J16").Value = 0.4325699745547074
$ws.Range("K16").Value = 0.08905852417302799
$ws.Range("M16").Value = 0.03307888040712468
$ws.Range("N16").Value = 0.002544529262086514
$ws.Range("O16").Value = 0.06361323155216285
$ws.Range("S16").Value = 0.07888040712468193
$ws.Range("F17").Value = 0.01769911504424779
$ws.Range("H17").Value = 0.2271386430678466
$ws.Range("I17").Value = 0.08554572271386431
$ws.Range("J17").Value = 0.4444444444444444
$ws.Range("K17").Value = 0.08062930186823992
$ws.Range("M17").Value = 0.01966568338249754
$ws.Range("N17").Value = 0.0009832841691248771
$ws.Range("O17").Value = 0.05703048180924287
$ws.Range("S17").Value = 0.06686332350049164
$ws.Range("F18").Value = 0.02586206896551724
$ws.Range("H18").Value = 0.1896551724137931
$ws.Range("I18").Value = 0.07327586206896551
$ws.Range("J18").Value = 0.478448275862069
$ws.Range("K18").Value = 0.09482758620689655
$ws.Range("M18").Value = 0.01939655172413793
$ws.Range("N18").Value = 0.002155172413793103
$ws.Range("O18").Value = 0.05818965517241379
$ws.Range("S18").Value = 0.05818965517241379
$ws.Range("F19").Value = 0.01843547583457898
$ws.Range("H19").Value = 0.2341803687095167
$ws.Range("I19").Value = 0.08320876930742402
$ws.Range("J19").Value = 0.3911310413552566
$ws.Range("K19").Value = 0.09466865969108122
$ws.Range("M19").Value = 0.0293971101145989
$ws.Range("N19").Value = 0.0009965122072745391
$ws.Range("O19").Value = 0.05630293971101146
$ws.Range("S19").Value = 0.0916791230692576
